$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
